# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..39 (replacing the old Strike# counts)
$kValues = @(7,9,6,8,9,4,9,9,6,11,5,6,7,7,6,9,8,4,6,12,8,14,13,12,7,8,6,9,6,6,4,10,4,5,4,3,4,2)

$row = 2
foreach ($val in $kValues) {
    $ws.Range("G$row").Value = $val
    $row++
}
